$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61
$ws.Cells.Item(61, 2).Value2 = 6905571
$ws.Cells.Item(61, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(61, 4).Value2 = 45130.72916666666
$ws.Cells.Item(61, 5).Value2 = "FBC Melgar"
$ws.Cells.Item(61, 6).Value2 = "Sporting Cristal"
$ws.Cells.Item(61, 7).Value2 = 1
$ws.Cells.Item(61, 8).Value2 = 1
$ws.Cells.Item(61, 9).Value2 = 1
$ws.Cells.Item(61, 10).Value2 = 0
$ws.Cells.Item(61, 11).Value2 = "D"
$ws.Cells.Item(61, 12).Value2 = 2.1
$ws.Cells.Item(61, 13).Value2 = 3.4
$ws.Cells.Item(61, 14).Value2 = 3
$ws.Cells.Item(61, 15).Value2 = 1.75
$ws.Cells.Item(61, 16).Value2 = 3.8
$ws.Cells.Item(61, 17).Value2 = 4.75
$ws.Cells.Item(61, 18).Value2 = -0.75
$ws.Cells.Item(61, 19).Value2 = 1.95
$ws.Cells.Item(61, 20).Value2 = 1.85
$ws.Cells.Item(61, 21).Value2 = 2.5
$ws.Cells.Item(61, 22).Value2 = 1.95
$ws.Cells.Item(61, 23).Value2 = 1.85
$ws.Cells.Item(61, 24).Value2 = -1
$ws.Cells.Item(61, 25).Value2 = 2.8
$ws.Cells.Item(61, 26).Value2 = -1
$ws.Cells.Item(61, 27).Value2 = -1
$ws.Cells.Item(61, 28).Value2 = 0.8500000000000001
$ws.Cells.Item(61, 29).Value2 = -1
$ws.Cells.Item(61, 30).Value2 = 0.8500000000000001

# Row 62
$ws.Cells.Item(62, 2).Value2 = 6905578
$ws.Cells.Item(62, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(62, 4).Value2 = 45130.72916666666
$ws.Cells.Item(62, 5).Value2 = "AD Tarma"
$ws.Cells.Item(62, 6).Value2 = "Atletico Grau"
$ws.Cells.Item(62, 7).Value2 = 1
$ws.Cells.Item(62, 8).Value2 = 0
$ws.Cells.Item(62, 9).Value2 = 1
$ws.Cells.Item(62, 10).Value2 = 0
$ws.Cells.Item(62, 11).Value2 = "H"
$ws.Cells.Item(62, 12).Value2 = 1.75
$ws.Cells.Item(62, 13).Value2 = 3.6
$ws.Cells.Item(62, 14).Value2 = 4
$ws.Cells.Item(62, 15).Value2 = 1.571
$ws.Cells.Item(62, 16).Value2 = 4.2
$ws.Cells.Item(62, 17).Value2 = 5.75
$ws.Cells.Item(62, 18).Value2 = -1
$ws.Cells.Item(62, 19).Value2 = 1.975
$ws.Cells.Item(62, 20).Value2 = 1.825
$ws.Cells.Item(62, 21).Value2 = 2.5
$ws.Cells.Item(62, 22).Value2 = 1.8
$ws.Cells.Item(62, 23).Value2 = 2
$ws.Cells.Item(62, 24).Value2 = 0.571
$ws.Cells.Item(62, 25).Value2 = -1
$ws.Cells.Item(62, 26).Value2 = -1
$ws.Cells.Item(62, 27).Value2 = 0
$ws.Cells.Item(62, 28).Value2 = 0
$ws.Cells.Item(62, 29).Value2 = -1
$ws.Cells.Item(62, 30).Value2 = 1

# Row 184
$ws.Cells.Item(184, 2).Value2 = 7384625
$ws.Cells.Item(184, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(184, 4).Value2 = 45228.70833333334
$ws.Cells.Item(184, 5).Value2 = "AD Tarma"
$ws.Cells.Item(184, 6).Value2 = "Carlos Manucci"
$ws.Cells.Item(184, 7).Value2 = 0
$ws.Cells.Item(184, 8).Value2 = 0
$ws.Cells.Item(184, 9).Value2 = 0
$ws.Cells.Item(184, 10).Value2 = 0
$ws.Cells.Item(184, 11).Value2 = "D"
$ws.Cells.Item(184, 12).Value2 = 1.5
$ws.Cells.Item(184, 13).Value2 = 3.75
$ws.Cells.Item(184, 14).Value2 = 7
$ws.Cells.Item(184, 15).Value2 = 1.363
$ws.Cells.Item(184, 16).Value2 = 4.333
$ws.Cells.Item(184, 17).Value2 = 9.5
$ws.Cells.Item(184, 18).Value2 = -1.25
$ws.Cells.Item(184, 19).Value2 = 1.875
$ws.Cells.Item(184, 20).Value2 = 1.925
$ws.Cells.Item(184, 21).Value2 = 2.5
$ws.Cells.Item(184, 22).Value2 = 1.8
$ws.Cells.Item(184, 23).Value2 = 2
$ws.Cells.Item(184, 24).Value2 = -1
$ws.Cells.Item(184, 25).Value2 = 3.333
$ws.Cells.Item(184, 26).Value2 = -1
$ws.Cells.Item(184, 27).Value2 = -1
$ws.Cells.Item(184, 28).Value2 = 0.925
$ws.Cells.Item(184, 29).Value2 = -1
$ws.Cells.Item(184, 30).Value2 = 1

# Row 185
$ws.Cells.Item(185, 2).Value2 = 7384628
$ws.Cells.Item(185, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(185, 4).Value2 = 45228.70833333334
$ws.Cells.Item(185, 5).Value2 = "Deportivo Binacional"
$ws.Cells.Item(185, 6).Value2 = "FBC Melgar"
$ws.Cells.Item(185, 7).Value2 = 1
$ws.Cells.Item(185, 8).Value2 = 2
$ws.Cells.Item(185, 9).Value2 = 1
$ws.Cells.Item(185, 10).Value2 = 1
$ws.Cells.Item(185, 11).Value2 = "A"
$ws.Cells.Item(185, 12).Value2 = 2.75
$ws.Cells.Item(185, 13).Value2 = 3.3
$ws.Cells.Item(185, 14).Value2 = 2.375
$ws.Cells.Item(185, 15).Value2 = 3.3
$ws.Cells.Item(185, 16).Value2 = 3.6
$ws.Cells.Item(185, 17).Value2 = 2
$ws.Cells.Item(185, 18).Value2 = 0.5
$ws.Cells.Item(185, 19).Value2 = 1.8
$ws.Cells.Item(185, 20).Value2 = 2
$ws.Cells.Item(185, 21).Value2 = 2.75
$ws.Cells.Item(185, 22).Value2 = 1.975
$ws.Cells.Item(185, 23).Value2 = 1.875
$ws.Cells.Item(185, 24).Value2 = -1
$ws.Cells.Item(185, 25).Value2 = -1
$ws.Cells.Item(185, 26).Value2 = 1
$ws.Cells.Item(185, 27).Value2 = -1
$ws.Cells.Item(185, 28).Value2 = 1
$ws.Cells.Item(185, 29).Value2 = 0.4875
$ws.Cells.Item(185, 30).Value2 = -0.5

# Row 186
$ws.Cells.Item(186, 2).Value2 = 7384627
$ws.Cells.Item(186, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(186, 4).Value2 = 45228.70833333334
$ws.Cells.Item(186, 5).Value2 = "Universitario de Deportes"
$ws.Cells.Item(186, 6).Value2 = "Sport Huancayo"
$ws.Cells.Item(186, 7).Value2 = 2
$ws.Cells.Item(186, 8).Value2 = 0
$ws.Cells.Item(186, 9).Value2 = 1
$ws.Cells.Item(186, 10).Value2 = 0
$ws.Cells.Item(186, 11).Value2 = "H"
$ws.Cells.Item(186, 12).Value2 = 1.25
$ws.Cells.Item(186, 13).Value2 = 5
$ws.Cells.Item(186, 14).Value2 = 12
$ws.Cells.Item(186, 15).Value2 = 1.181
$ws.Cells.Item(186, 16).Value2 = 6
$ws.Cells.Item(186, 17).Value2 = 13
$ws.Cells.Item(186, 18).Value2 = -1.75
$ws.Cells.Item(186, 19).Value2 = 1.8
$ws.Cells.Item(186, 20).Value2 = 2
$ws.Cells.Item(186, 21).Value2 = 2.75
$ws.Cells.Item(186, 22).Value2 = 1.85
$ws.Cells.Item(186, 23).Value2 = 1.95
$ws.Cells.Item(186, 24).Value2 = 0.181
$ws.Cells.Item(186, 25).Value2 = -1
$ws.Cells.Item(186, 26).Value2 = -1
$ws.Cells.Item(186, 27).Value2 = 0.4
$ws.Cells.Item(186, 28).Value2 = -0.5
$ws.Cells.Item(186, 29).Value2 = -1
$ws.Cells.Item(186, 30).Value2 = 0.95

# Row 187
$ws.Cells.Item(187, 2).Value2 = 7384626
$ws.Cells.Item(187, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(187, 4).Value2 = 45228.70833333334
$ws.Cells.Item(187, 5).Value2 = "Sporting Cristal"
$ws.Cells.Item(187, 6).Value2 = "Alianza Atletico"
$ws.Cells.Item(187, 7).Value2 = 3
$ws.Cells.Item(187, 8).Value2 = 0
$ws.Cells.Item(187, 9).Value2 = 3
$ws.Cells.Item(187, 10).Value2 = 0
$ws.Cells.Item(187, 11).Value2 = "H"
$ws.Cells.Item(187, 12).Value2 = 1.3
$ws.Cells.Item(187, 13).Value2 = 5
$ws.Cells.Item(187, 14).Value2 = 9
$ws.Cells.Item(187, 15).Value2 = 1.166
$ws.Cells.Item(187, 16).Value2 = 6.5
$ws.Cells.Item(187, 17).Value2 = 13
$ws.Cells.Item(187, 18).Value2 = -2
$ws.Cells.Item(187, 19).Value2 = 1.85
$ws.Cells.Item(187, 20).Value2 = 1.95
$ws.Cells.Item(187, 21).Value2 = 3.25
$ws.Cells.Item(187, 22).Value2 = 2
$ws.Cells.Item(187, 23).Value2 = 1.8
$ws.Cells.Item(187, 24).Value2 = 0.1659999999999999
$ws.Cells.Item(187, 25).Value2 = -1
$ws.Cells.Item(187, 26).Value2 = -1
$ws.Cells.Item(187, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(187, 28).Value2 = -1
$ws.Cells.Item(187, 29).Value2 = -0.5
$ws.Cells.Item(187, 30).Value2 = 0.4

# Row 188
$ws.Cells.Item(188, 2).Value2 = 7384629
$ws.Cells.Item(188, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(188, 4).Value2 = 45228.70833333334
$ws.Cells.Item(188, 5).Value2 = "Deportivo Garcilaso"
$ws.Cells.Item(188, 6).Value2 = "Alianza Lima"
$ws.Cells.Item(188, 7).Value2 = 0
$ws.Cells.Item(188, 8).Value2 = 1
$ws.Cells.Item(188, 9).Value2 = 0
$ws.Cells.Item(188, 10).Value2 = 1
$ws.Cells.Item(188, 11).Value2 = "A"
$ws.Cells.Item(188, 12).Value2 = 2.625
$ws.Cells.Item(188, 13).Value2 = 3.3
$ws.Cells.Item(188, 14).Value2 = 2.5
$ws.Cells.Item(188, 15).Value2 = 2.7
$ws.Cells.Item(188, 16).Value2 = 3.4
$ws.Cells.Item(188, 17).Value2 = 2.375
$ws.Cells.Item(188, 18).Value2 = 0
$ws.Cells.Item(188, 19).Value2 = 2.025
$ws.Cells.Item(188, 20).Value2 = 1.775
$ws.Cells.Item(188, 21).Value2 = 2.25
$ws.Cells.Item(188, 22).Value2 = 1.825
$ws.Cells.Item(188, 23).Value2 = 1.975
$ws.Cells.Item(188, 24).Value2 = -1
$ws.Cells.Item(188, 25).Value2 = -1
$ws.Cells.Item(188, 26).Value2 = 1.375
$ws.Cells.Item(188, 27).Value2 = -1
$ws.Cells.Item(188, 28).Value2 = 0.7749999999999999
$ws.Cells.Item(188, 29).Value2 = -1
$ws.Cells.Item(188, 30).Value2 = 0.9750000000000001

# Row 228
$ws.Cells.Item(228, 2).Value2 = 7818816
$ws.Cells.Item(228, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(228, 4).Value2 = 45346.70833333334
$ws.Cells.Item(228, 5).Value2 = "UTC Cajamarca"
$ws.Cells.Item(228, 6).Value2 = "Universitario de Deportes"
$ws.Cells.Item(228, 7).Value2 = 0
$ws.Cells.Item(228, 8).Value2 = 0
$ws.Cells.Item(228, 9).Value2 = 0
$ws.Cells.Item(228, 10).Value2 = 0
$ws.Cells.Item(228, 11).Value2 = "D"
$ws.Cells.Item(228, 12).Value2 = 3.3
$ws.Cells.Item(228, 13).Value2 = 3.3
$ws.Cells.Item(228, 14).Value2 = 2.1
$ws.Cells.Item(228, 15).Value2 = 4.5
$ws.Cells.Item(228, 16).Value2 = 3.2
$ws.Cells.Item(228, 17).Value2 = 1.95
$ws.Cells.Item(228, 18).Value2 = 0.5
$ws.Cells.Item(228, 19).Value2 = 2
$ws.Cells.Item(228, 20).Value2 = 1.85
$ws.Cells.Item(228, 21).Value2 = 2
$ws.Cells.Item(228, 22).Value2 = 1.775
$ws.Cells.Item(228, 23).Value2 = 2.1
$ws.Cells.Item(228, 24).Value2 = -1
$ws.Cells.Item(228, 25).Value2 = 2.2
$ws.Cells.Item(228, 26).Value2 = -1
$ws.Cells.Item(228, 27).Value2 = 1
$ws.Cells.Item(228, 28).Value2 = -1
$ws.Cells.Item(228, 29).Value2 = -1
$ws.Cells.Item(228, 30).Value2 = 1.1

# Row 229
$ws.Cells.Item(229, 2).Value2 = 7818817
$ws.Cells.Item(229, 3).Value2 = "Peru Liga 1"
$ws.Cells.Item(229, 4).Value2 = 45346.70833333334
$ws.Cells.Item(229, 5).Value2 = "Sport Boys"
$ws.Cells.Item(229, 6).Value2 = "Cusco FC"
$ws.Cells.Item(229, 7).Value2 = 3
$ws.Cells.Item(229, 8).Value2 = 0
$ws.Cells.Item(229, 9).Value2 = 2
$ws.Cells.Item(229, 10).Value2 = 0
$ws.Cells.Item(229, 11).Value2 = "H"
$ws.Cells.Item(229, 12).Value2 = 2.2
$ws.Cells.Item(229, 13).Value2 = 3.2
$ws.Cells.Item(229, 14).Value2 = 3.2
$ws.Cells.Item(229, 15).Value2 = 1.6
$ws.Cells.Item(229, 16).Value2 = 3.75
$ws.Cells.Item(229, 17).Value2 = 5.75
$ws.Cells.Item(229, 18).Value2 = -0.75
$ws.Cells.Item(229, 19).Value2 = 1.85
$ws.Cells.Item(229, 20).Value2 = 2
$ws.Cells.Item(229, 21).Value2 = 2.5
$ws.Cells.Item(229, 22).Value2 = 1.975
$ws.Cells.Item(229, 23).Value2 = 1.875
$ws.Cells.Item(229, 24).Value2 = 0.6000000000000001
$ws.Cells.Item(229, 25).Value2 = -1
$ws.Cells.Item(229, 26).Value2 = -1
$ws.Cells.Item(229, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(229, 28).Value2 = -1
$ws.Cells.Item(229, 29).Value2 = 0.9750000000000001
$ws.Cells.Item(229, 30).Value2 = 1.1
